$wb = $excel.ActiveWorkbook

# "Overview" sheet: refresh the rolled-up "Latest HO Xliff Generate Date" for the
# first file (1dfcb511-...) to reflect the newly generated handback.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-09-07 03:01:33"

# "zh-cn" handback report: refresh handoff/handback datetimes for the first file.
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-09-07 03:01:28"
$zhcn.Range("K2").Value = "2016-09-07 03:01:47"

# "de-de" handback report: refresh handoff/handback datetimes for the first file.
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-09-07 03:01:33"
$dede.Range("K2").Value = "2016-09-07 03:01:56"
